$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet has an Excel Table ("Tabela1") covering A1:D70.
# Add a new row to the table for "Aula 68 - Padronizando as URLs internas".
$tbl = $ws.ListObjects.Item(1)
$newListRow = $tbl.ListRows.Add()
$r = $newListRow.Range.Row

# Copy the formatting of the previous last data row (row 70) onto the new row
# so that fonts / wrap / etc. match the rest of the table.
$ws.Range("A70:D70").Copy()
$ws.Range("A71:D71").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the new row's values (observação, nome da aula, sessão, then aula,
# mirroring the order the shared strings were originally authored in so the
# shared-string table is built up in the same sequence).
$obsText = "1:34`r`nsubstituido todas links e tags href e src que apontam para as paginas HTML da aplicação pelo padrão do thymeleaf:`r`nth:href=""@{/suapaginaaqui}`r`nth:src=""@{/seucaminhoaqui}"
$ws.Cells.Item($r, 4).Value = $obsText
$ws.Cells.Item($r, 3).Value = "`n68. Padronizando as URLs internas"
$ws.Cells.Item($r, 2).Value = "14. Final"
$ws.Cells.Item($r, 1).Value = 68

# Match the row height used for similarly-sized observation text elsewhere in the sheet.
$ws.Rows.Item($r).RowHeight = 90
